$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('B2').Value = 'Bitcoin'
$ws.Range('D2').Value = '25.757.00'
$ws.Range('E2').Value = '  +0.24%  '

# Row 3
$ws.Range('B3').Value = 'Ethereum'
$ws.Range('D3').Value = '1.748.78'
$ws.Range('E3').Value = '  +0.48%  '

# Row 4
$ws.Range('B4').Value = 'TetherUSD'
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.002'
$ws.Range('E4').Value = '  +0.10%  '

# Row 5
$ws.Range('B5').Value = 'BNB'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '236.65'
$ws.Range('E5').Value = '  -0.92%  '

# Row 6
$ws.Range('B6').Value = 'USDC'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.001'
$ws.Range('E6').Value = '  -0.02%  '

# Row 7
$ws.Range('B7').Value = 'XRP'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5054'
$ws.Range('E7').Value = '  +2.98%  '

# Row 8
$ws.Range('B8').Value = 'Cardano'
$ws.Range('C8').Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2604'
$ws.Range('E8').Value = '  +8.03%  '

# Row 9
$ws.Range('B9').Value = 'Dogecoin'
$ws.Range('C9').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06173'
$ws.Range('E9').Value = '  +3.04%  '

# Row 10
$ws.Range('B10').Value = 'WrappedEther'
$ws.Range('C10').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D10').Value = '1.751.73'
$ws.Range('E10').Value = '  +0.61%  '

# Row 11
$ws.Range('B11').Value = 'TRON'
$ws.Range('C11').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.06931'
$ws.Range('E11').Value = '  +2.05%  '

# Row 12
$ws.Range('B12').Value = 'Solana'
$ws.Range('C12').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '15.36'
$ws.Range('E12').Value = '  +3.91%  '

# Row 13
$ws.Range('B13').Value = 'Polygon'
$ws.Range('C13').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.6053'
$ws.Range('E13').Value = '  +4.33%  '

# Row 14
$ws.Range('B14').Value = 'Litecoin'
$ws.Range('C14').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '78.40'
$ws.Range('E14').Value = '  +2.41%  '

# Row 15
$ws.Range('B15').Value = 'Polkadot'
$ws.Range('C15').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '4.452'
$ws.Range('E15').Value = '  +0.31%  '

# Row 16
$ws.Range('B16').Value = 'BinanceUSD'
$ws.Range('C16').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.002'
$ws.Range('E16').Value = '  +0.18%  '

# Row 17
$ws.Range('B17').Value = 'Dai'
$ws.Range('C17').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.001'
$ws.Range('E17').Value = '  -0.02%  '

# Row 18
$ws.Range('B18').Value = 'WrappedBTC'
$ws.Range('C18').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D18').Value = '25.781.48'
$ws.Range('E18').Value = '  +0.20%  '

# Row 19
$ws.Range('B19').Value = 'Avalanche'
$ws.Range('C19').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '11.64'
$ws.Range('E19').Value = '  +1.48%  '

# Row 20
$ws.Range('B20').Value = 'ShibaInu'
$ws.Range('C20').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.000006691'
$ws.Range('E20').Value = '  +4.75%  '

# Row 21
$ws.Range('B21').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C21').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D21').Value = '1.976.06'
$ws.Range('E21').Value = '  +0.87%  '

# Row 22
$ws.Range('B22').Value = 'Uniswap'
$ws.Range('C22').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.048'
$ws.Range('E22').Value = '  +2.83%  '

# Row 23
$ws.Range('B23').Value = 'Cosmos'
$ws.Range('C23').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '8.189'
$ws.Range('E23').Value = '  +4.42%  '

# Row 24
$ws.Range('B24').Value = 'Chainlink'
$ws.Range('C24').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '5.159'
$ws.Range('E24').Value = '  +1.34%  '

# Row 25
$ws.Range('B25').Value = 'Monero'
$ws.Range('C25').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '137.68'
$ws.Range('E25').Value = '  +1.14%  '

# Row 26
$ws.Range('B26').Value = 'Toncoin'
$ws.Range('C26').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.458'
$ws.Range('E26').Value = '  +0.10%  '

# Row 27
$ws.Range('B27').Value = 'EthereumClassic'
$ws.Range('C27').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '15.09'
$ws.Range('E27').Value = '  +4.41%  '

# Row 28
$ws.Range('B28').Value = 'LidoDAOToken'
$ws.Range('C28').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.802'
$ws.Range('E28').Value = '  -2.31%  '

# Row 29
$ws.Range('B29').Value = 'BitcoinCash'
$ws.Range('C29').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '102.23'
$ws.Range('E29').Value = '  +2.40%  '

# Row 30
$ws.Range('B30').Value = 'Stellar'
$ws.Range('C30').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.08270'
$ws.Range('E30').Value = '  +2.15%  '

# Row 31
$ws.Range('B31').Value = 'InternetComputer(DFINITY)'
$ws.Range('C31').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.697'
$ws.Range('E31').Value = '  -0.65%  '

# Row 32
$ws.Range('B32').Value = 'Filecoin'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.394'
$ws.Range('E32').Value = '  +0.73%  '

# Row 33
$ws.Range('B33').Value = 'Hedera'
$ws.Range('C33').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.04357'
$ws.Range('E33').Value = '  -0.09%  '

# Row 34
$ws.Range('B34').Value = 'Frax'
$ws.Range('C34').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.0000'
$ws.Range('E34').Value = '  +0.00%  '

# Row 35
$ws.Range('B35').Value = 'HuobiToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.653'
$ws.Range('E35').Value = '  -1.50%  '

# Row 36
$ws.Range('B36').Value = 'ARBITRUM'
$ws.Range('C36').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.9969'
$ws.Range('E36').Value = '  -2.01%  '

# Row 37
$ws.Range('B37').Value = 'ImmutableX'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.6009'
$ws.Range('E37').Value = '  -0.01%  '

# Row 38
$ws.Range('B38').Value = 'MXToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.699'
$ws.Range('E38').Value = '  -1.19%  '

# Row 39
$ws.Range('B39').Value = 'RenderToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.957'
$ws.Range('E39').Value = '  -4.77%  '

# Row 40
$ws.Range('B40').Value = 'VeChain'
$ws.Range('C40').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.01549'
$ws.Range('E40').Value = '  +4.38%  '

# Row 41
$ws.Range('B41').Value = 'PaxDollar'
$ws.Range('C41').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.001'
$ws.Range('E41').Value = '  +0.00%  '

# Row 42
$ws.Range('B42').Value = 'PaxosStandard'
$ws.Range('C42').Value = 'https://coinranking.com/coin/B8xT718SbVhhh+paxosstandard-pax'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.002'
$ws.Range('E42').Value = '  +0.08%  '

# Row 43
$ws.Range('B43').Value = 'Quant'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '102.69'
$ws.Range('E43').Value = '  -0.44%  '

# Row 44
$ws.Range('B44').Value = 'TrustWalletToken'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.7572'
$ws.Range('E44').Value = '  -3.48%  '

# Row 45
$ws.Range('B45').Value = 'TheSandbox'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.3803'
$ws.Range('E45').Value = '  +0.14%  '

# Row 46
$ws.Range('B46').Value = 'FraxShare'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '4.845'
$ws.Range('E46').Value = '  -5.65%  '

# Row 47
$ws.Range('B47').Value = 'Cronos'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.05495'
$ws.Range('E47').Value = '  +7.92%  '

# Row 48
$ws.Range('B48').Value = 'Algorand'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.1079'
$ws.Range('E48').Value = '  +1.46%  '

# Row 49
$ws.Range('B49').Value = 'Elrond'
$ws.Range('C49').Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '30.13'
$ws.Range('E49').Value = '  -0.29%  '

# Row 50
$ws.Range('B50').Value = 'Aptos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '5.915'
$ws.Range('E50').Value = '  -1.39%  '

# Row 51
$ws.Range('B51').Value = 'USDD'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.002'
$ws.Range('E51').Value = '  +0.16%  '
